# Insert a new weekly price record as row 58 (Macroferia Regional de Talca,
# Arándano (blue)), pushing the existing rows 58-92 down to 59-93.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(58).Insert()

$ws.Range("A58").Value = 5
$ws.Range("B58").Value = "Macroferia Regional de Talca"
$ws.Range("C58").Value = "Maule"
$ws.Range("D58").Value = 44907
$ws.Range("E58").Value = 7
$ws.Range("F58").Value = "Fruta"
$ws.Range("G58").Value = 100101
$ws.Range("H58").Value = "Berries"
$ws.Range("I58").Value = 100101001
$ws.Range("J58").Value = "Arándano (blue)"
$ws.Range("K58").Value = "Sin especificar"
$ws.Range("L58").Value = "Primera"
$ws.Range("M58").Value = 180
$ws.Range("N58").Value = 3000
$ws.Range("O58").Value = 3000
$ws.Range("P58").Value = 3000
$ws.Range("Q58").Value = "$/bandeja 2 kilos"
$ws.Range("R58").Value = "Provincia de Curicó"
$ws.Range("S58").Value = 1500
$ws.Range("T58").Value = 2
